# Fill in the "P3 # Training" (column E) and "P3 # Testing" (column L) counts
# for rows 2-51, fix a typo in one of the shared-string labels, and move the
# active-cell selection, matching the "P3 Training and Testing Added" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E = "P3 # Training" values for rows 2..51
$p3Training = @(24,24,24,24,35,29,27,24,26,25,26,28,28,26,25,25,30,25,26,24,28,24,25,27,26,25,57,56,26,27,25,24,25,25,24,24,24,24,24,24,24,24,24,24,24,24,24,26,25,26)

# Column L = "P3 # Testing" values for rows 2..51
$p3Testing = @(6,6,6,6,12,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6,6)

for ($i = 0; $i -lt $p3Training.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $p3Training[$i]
    $ws.Cells.Item($row, 12).Value = $p3Testing[$i]
}

# Fix typo in the shared-string note
$ws.Range("P4").Value = "*P6 has 1 diff signer"

# Move the active-cell selection from I7 to I6
$ws.Range("I6").Select()
